# Applies the "Updated cryptos list" data refresh to the active worksheet.
# Updates Price (D) and Volume(1h) (E) columns for each coin row, and swaps the
# RocketPoolETH / BitcoinSV rows (48-49) to reflect the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes $value into (row,col) and keeps it stored as text, since many
# of the values (prices such as "43.387.19" or "97.62") would otherwise be
# auto-converted to numbers by Excel. Forcing the NumberFormat to "@" before
# the assignment keeps the literal string, and ClearFormats afterwards removes
# the temporary text format so the cell keeps its original (default) style.
function Set-TextValue($row, $col, $value) {
    $rng = $ws.Cells.Item($row, $col)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
Set-TextValue 2 4 '43.387.19'
Set-TextValue 2 5 '  +1.41%  '
# Row 3
Set-TextValue 3 4 '2.602.31'
Set-TextValue 3 5 '  +3.12%  '
# Row 4
Set-TextValue 4 5 '  +0.10%  '
# Row 5
Set-TextValue 5 4 '316.69'
Set-TextValue 5 5 '  -0.06%  '
# Row 6
Set-TextValue 6 4 '97.62'
Set-TextValue 6 5 '  +3.06%  '
# Row 7
Set-TextValue 7 4 '0.578'
Set-TextValue 7 5 '  -0.28%  '
# Row 8
Set-TextValue 8 5 '  +0.02%  '
# Row 9
Set-TextValue 9 4 '0.543'
Set-TextValue 9 5 '  +2.43%  '
# Row 10
Set-TextValue 10 4 '35.96'
Set-TextValue 10 5 '  +0.23%  '
# Row 11
Set-TextValue 11 4 '0.0816'
Set-TextValue 11 5 '  +0.82%  '
# Row 12
Set-TextValue 12 4 '7.55'
Set-TextValue 12 5 '  -0.09%  '
# Row 13
Set-TextValue 13 4 '3.002.73'
Set-TextValue 13 5 '  +3.11%  '
# Row 14
Set-TextValue 14 5 '  -1.12%  '
# Row 15
Set-TextValue 15 4 '2.627.16'
Set-TextValue 15 5 '  +3.86%  '
# Row 16
Set-TextValue 16 4 '15.31'
Set-TextValue 16 5 '  +0.64%  '
# Row 17
Set-TextValue 17 4 '0.850'
Set-TextValue 17 5 '  +0.23%  '
# Row 18
Set-TextValue 18 4 '43.534.84'
Set-TextValue 18 5 '  +1.55%  '
# Row 19
Set-TextValue 19 4 '6.88'
Set-TextValue 19 5 '  +2.80%  '
# Row 20
Set-TextValue 20 5 '  -2.13%  '
# Row 21
Set-TextValue 21 4 '0.0₃0970'
Set-TextValue 21 5 '  +0.64%  '
# Row 22
Set-TextValue 22 4 '69.89'
Set-TextValue 22 5 '  +0.14%  '
# Row 23
Set-TextValue 23 4 '255.45'
Set-TextValue 23 5 '  +1.71%  '
# Row 24
Set-TextValue 24 5 '  +1.22%  '
# Row 25
Set-TextValue 25 5 '  +3.37%  '
# Row 26
Set-TextValue 26 4 '27.53'
Set-TextValue 26 5 '  +2.88%  '
# Row 27
Set-TextValue 27 5 '  -0.21%  '
# Row 28
Set-TextValue 28 5 '  +0.65%  '
# Row 29
Set-TextValue 29 4 '40.65'
Set-TextValue 29 5 '  +0.57%  '
# Row 30
Set-TextValue 30 4 '10.33'
Set-TextValue 30 5 '  +0.49%  '
# Row 31
Set-TextValue 31 4 '5.92'
Set-TextValue 31 5 '  -1.15%  '
# Row 32
Set-TextValue 32 4 '157.59'
Set-TextValue 32 5 '  +0.51%  '
# Row 33
Set-TextValue 33 4 '3.51'
Set-TextValue 33 5 '  +6.94%  '
# Row 34
Set-TextValue 34 4 '2.17'
Set-TextValue 34 5 '  +2.74%  '
# Row 35
Set-TextValue 35 4 '0.0813'
Set-TextValue 35 5 '  +3.22%  '
# Row 36
Set-TextValue 36 5 '  +3.63%  '
# Row 37
Set-TextValue 37 4 '18.87'
Set-TextValue 37 5 '  -0.37%  '
# Row 38
Set-TextValue 38 4 '2.53'
Set-TextValue 38 5 '  +10.28%  '
# Row 39
Set-TextValue 39 5 '  +0.36%  '
# Row 40
Set-TextValue 40 5 '  +0.08%  '
# Row 41
Set-TextValue 41 4 '23.02'
Set-TextValue 41 5 '  -2.38%  '
# Row 42
Set-TextValue 42 4 '4.04'
Set-TextValue 42 5 '  +7.42%  '
# Row 43
Set-TextValue 43 4 '0.0305'
Set-TextValue 43 5 '  +0.25%  '
# Row 44
Set-TextValue 44 5 '  +0.08%  '
# Row 45
Set-TextValue 45 4 '3.28'
Set-TextValue 45 5 '  -0.71%  '
# Row 46
Set-TextValue 46 4 '2.022.52'
Set-TextValue 46 5 '  +0.07%  '
# Row 47
Set-TextValue 47 4 '9.02'
Set-TextValue 47 5 '  +2.74%  '
# Row 48
Set-TextValue 48 2 'RocketPoolETH'
Set-TextValue 48 3 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 48 4 '2.856.11'
Set-TextValue 48 5 '  +3.22%  '
# Row 49
Set-TextValue 49 2 'BitcoinSV'
Set-TextValue 49 3 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue 49 4 '83.96'
Set-TextValue 49 5 '  -1.94%  '
# Row 50
Set-TextValue 50 4 '75.25'
Set-TextValue 50 5 '  +2.20%  '
# Row 51
Set-TextValue 51 4 '0.194'
Set-TextValue 51 5 '  +2.75%  '
